# ORG_CIRCLE.xlsx — add three new tracking columns (D:F) to Sheet1's header
# row: ORG_CIR_IDENOLD, ORG_CIR_IDENNEW, ORG_CIR_STATUS. These land as new
# shared-string entries and extend the sheet's used range from A1:C6 to
# A1:F6 (only the header row gets values; the data rows 2-6 stay as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "ORG_CIR_IDENOLD"
$ws.Range("E1").Value = "ORG_CIR_IDENNEW"
$ws.Range("F1").Value = "ORG_CIR_STATUS"

# Leave the cursor on F9, matching the saved workbook view.
$ws.Range("F9").Select()
